$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aguilar Prototype")

# GDP Nowcast forecast value update
$ws.Range("F7").Value = 0.3747070985527221

# Row 29 - 5yr, 5yr Forward (T5YIFR)
$ws.Range("N29").Value = "2025-11-04"
$ws.Range("Q29").Value = 2.2

# Row 30 - 10yr TIPS (T10YIE)
$ws.Range("N30").Value = "2025-11-04"
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.31
$ws.Range("S30").Value = $null
$ws.Range("T30").Value = $null
$ws.Range("U30").Value = 2.3

# Row 47 - FFR (DFF)
$ws.Range("N47").Value = "2025-11-03"
$ws.Range("Q47").Value = 3.87
$ws.Range("R47").Value = 3.86
$ws.Range("S47").Value = 3.86
$ws.Range("T47").Value = 3.86
$ws.Range("U47").Value = 3.87

# Row 48 - 2y UST (DGS2)
$ws.Range("N48").Value = "2025-11-03"
$ws.Range("Q48").Value = 3.6
$ws.Range("R48").Value = $null
$ws.Range("S48").Value = $null
$ws.Range("T48").Value = 3.6
$ws.Range("U48").Value = 3.61

# Row 49 - 5y UST (DGS5)
$ws.Range("N49").Value = "2025-11-03"
$ws.Range("Q49").Value = 3.72
$ws.Range("R49").Value = $null
$ws.Range("S49").Value = $null
$ws.Range("T49").Value = 3.71
$ws.Range("U49").Value = 3.72

# Row 50 - 10y UST (DGS10)
$ws.Range("N50").Value = "2025-11-03"
$ws.Range("Q50").Value = 4.13
$ws.Range("R50").Value = $null
$ws.Range("S50").Value = $null
$ws.Range("T50").Value = 4.11
$ws.Range("U50").Value = 4.11

# Row 52 - BAA (DBAA)
$ws.Range("N52").Value = "2025-11-03"
$ws.Range("Q52").Value = 5.84
$ws.Range("R52").Value = $null
$ws.Range("S52").Value = $null
$ws.Range("T52").Value = 5.8
$ws.Range("U52").Value = 5.75
